# Add SNS support for new message.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Remove the broken external reference to OrderFlowersChatbot.xlsx
# ---------------------------------------------------------------------
$wb.BreakLink("OrderFlowersChatbot.xlsx", 1) | Out-Null

# ---------------------------------------------------------------------
# 2. OrderFlowersIntend (sheet3): insert an "Email" row above the
#    "slots" row, fix the stale external data-validation formulas.
# ---------------------------------------------------------------------
$wsOrderIntend = $wb.Worksheets.Item("OrderFlowersIntend")
$wsOrderIntend.Rows.Item(6).Insert() | Out-Null
$wsOrderIntend.Range("A6").Value = "Email"

$wsOrderIntend.Range("D8:D17").Validation.Delete() | Out-Null
$wsOrderIntend.Range("E8:E17").Validation.Delete() | Out-Null
$wsOrderIntend.Range("D9:D11").Validation.Add(3, 1, 1, "=Option!`$B`$2:`$B`$93") | Out-Null
$wsOrderIntend.Range("E9:E11").Validation.Add(3, 1, 1, "=Option!`$C`$2:`$C`$3") | Out-Null

$wsOrderIntend.Range("C20").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. MakeAppointmentIntend (sheet4): insert an "Email" row with an
#    actual mailto hyperlink to the bot author.
# ---------------------------------------------------------------------
$wsMakeAppt = $wb.Worksheets.Item("MakeAppointmentIntend")
$wsMakeAppt.Rows.Item(6).Insert() | Out-Null
$wsMakeAppt.Range("A6").Value = "Email"
$wsMakeAppt.Hyperlinks.Add($wsMakeAppt.Range("B6"), "mailto:cywong@vtc.edu.hk", [Type]::Missing, [Type]::Missing, "cywong@vtc.edu.hk") | Out-Null
$wsMakeAppt.Range("B6").Value = "cywong@vtc.edu.hk`ncy.gdoc@gmail.com"
$wsMakeAppt.Rows.Item(6).RowHeight = 30

$wsMakeAppt.Range("D9:D11").Validation.Delete() | Out-Null
$wsMakeAppt.Range("E9:E11").Validation.Delete() | Out-Null
$wsMakeAppt.Range("D9:D11").Validation.Add(3, 1, 1, "=Option!`$B`$2:`$B`$93") | Out-Null
$wsMakeAppt.Range("E9:E11").Validation.Add(3, 1, 1, "=Option!`$C`$2:`$C`$3") | Out-Null

$wsMakeAppt.Range("B6").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. FlowerTypes (sheet5): drop the stray H17 cell, fix the stale
#    external data-validation formula, move the selection.
# ---------------------------------------------------------------------
$wsFlowerTypes = $wb.Worksheets.Item("FlowerTypes")
$wsFlowerTypes.Range("H17").ClearContents() | Out-Null
$wsFlowerTypes.Range("B2").Validation.Delete() | Out-Null
$wsFlowerTypes.Range("B2").Validation.Add(3, 1, 1, "=Option!`$A`$2:`$A`$3") | Out-Null
$wsFlowerTypes.Range("B2").Select() | Out-Null

# ---------------------------------------------------------------------
# 5. AppointmentTypes (sheet6): drop the stray H17 cell, move the
#    selection (its data validation formula already points locally).
# ---------------------------------------------------------------------
$wsApptTypes = $wb.Worksheets.Item("AppointmentTypes")
$wsApptTypes.Range("H17").ClearContents() | Out-Null
$wsApptTypes.Range("B2").Select() | Out-Null

# ---------------------------------------------------------------------
# 6. OrderFlowersBot / ScheduleAppointmentBot: only the selection moves
#    (shared-string text shifts automatically once "Bot" is dropped).
# ---------------------------------------------------------------------
$wsOrderBot = $wb.Worksheets.Item("OrderFlowersBot")
$wsOrderBot.Range("B1").Select() | Out-Null

$wsScheduleBot = $wb.Worksheets.Item("ScheduleAppointmentBot")
$wsScheduleBot.Range("B1").Select() | Out-Null

# ---------------------------------------------------------------------
# Restore the originally active sheet/tab so tabSelected stays put.
# ---------------------------------------------------------------------
$wsOrderIntend.Activate() | Out-Null
